$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace / extend the demo data (tag_id / folder_id) -----------------
$data = @(
    @("personal_data", "bevnat"),
    @("personal_data", "statpop"),
    @("sensible_data", "bevnat"),
    @("population",    "bevnat"),
    @("population",    "statpop"),
    @("societe",       "statpop")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value() = $data[$i][0]
    $ws.Cells.Item($row, 2).Value() = $data[$i][1]
}

# --- Resize the table / autofilter to cover the new rows -----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B7"))

# --- Widen column A a bit so the longer labels fit -----------------------
$ws.Columns.Item(1).ColumnWidth = 15.83

# --- Update the current selection (cosmetic, mirrors the saved view) -----
$ws.Range("B9").Select()
